$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rewrite data rows 2-7 (row 2-3 updated values; rows 4-7 new) per updated TPM calculations
# Row 2
$ws.Range("A2").Value2 = "ECs"
$ws.Range("B2").Value2 = "Cd86"
$ws.Range("C2").Value2 = "Ctla4"
$ws.Range("D2").Value2 = "MuSCs"
$ws.Range("E2").Value2 = 1
$ws.Range("F2").Value2 = 0.3333333333333333
$ws.Range("G2").Value2 = 0.020562
$ws.Range("H2").Value2 = 0.061686
$ws.Range("I2").Value2 = 0.0002281281878049052
$ws.Range("J2").Value2 = 0.0002281281878049052
$ws.Range("K2").Value2 = 1
$ws.Range("L2").Value2 = 0.3333333333333333
$ws.Range("M2").Value2 = 0.005070666666666667
$ws.Range("N2").Value2 = 0.015212
$ws.Range("O2").Value2 = 0.008865695667453655
$ws.Range("P2").Value2 = 0.008865695667453653
$ws.Range("Q2").Value2 = 0.000104263048
$ws.Range("R2").Value2 = 0.000938367432
$ws.Range("S2").Value2 = 0.000002022515086246002
$ws.Range("T2").Value2 = 0.000002022515086246002

# Row 3
$ws.Range("A3").Value2 = "ECs"
$ws.Range("B3").Value2 = "Cd86"
$ws.Range("C3").Value2 = "Ctla4"
$ws.Range("D3").Value2 = "Resolving-Mac"
$ws.Range("E3").Value2 = 1
$ws.Range("F3").Value2 = 0.3333333333333333
$ws.Range("G3").Value2 = 0.020562
$ws.Range("H3").Value2 = 0.061686
$ws.Range("I3").Value2 = 0.0002281281878049052
$ws.Range("J3").Value2 = 0.0002281281878049052
$ws.Range("K3").Value2 = 3
$ws.Range("L3").Value2 = 1
$ws.Range("M3").Value2 = 0.5668716666666667
$ws.Range("N3").Value2 = 1.700615
$ws.Range("O3").Value2 = 0.9911343043325463
$ws.Range("P3").Value2 = 0.9911343043325463
$ws.Range("Q3").Value2 = 0.01165601521
$ws.Range("R3").Value2 = 0.10490413689
$ws.Range("S3").Value2 = 0.0002261056727186593
$ws.Range("T3").Value2 = 0.0002261056727186593

# Row 4
$ws.Range("A4").Value2 = "MuSCs"
$ws.Range("B4").Value2 = "Cd86"
$ws.Range("C4").Value2 = "Ctla4"
$ws.Range("D4").Value2 = "MuSCs"
$ws.Range("E4").Value2 = 1
$ws.Range("F4").Value2 = 0.3333333333333333
$ws.Range("G4").Value2 = 0.01154533333333333
$ws.Range("H4").Value2 = 0.034636
$ws.Range("I4").Value2 = 0.0001280914293812323
$ws.Range("J4").Value2 = 0.0001280914293812324
$ws.Range("K4").Value2 = 1
$ws.Range("L4").Value2 = 0.3333333333333333
$ws.Range("M4").Value2 = 0.005070666666666667
$ws.Range("N4").Value2 = 0.015212
$ws.Range("O4").Value2 = 0.008865695667453655
$ws.Range("P4").Value2 = 0.008865695667453653
$ws.Range("Q4").Value2 = 0.00005854253688888889
$ws.Range("R4").Value2 = 0.000526882832
$ws.Range("S4").Value2 = 0.000001135619630503137
$ws.Range("T4").Value2 = 0.000001135619630503137

# Row 5
$ws.Range("A5").Value2 = "MuSCs"
$ws.Range("B5").Value2 = "Cd86"
$ws.Range("C5").Value2 = "Ctla4"
$ws.Range("D5").Value2 = "Resolving-Mac"
$ws.Range("E5").Value2 = 1
$ws.Range("F5").Value2 = 0.3333333333333333
$ws.Range("G5").Value2 = 0.01154533333333333
$ws.Range("H5").Value2 = 0.034636
$ws.Range("I5").Value2 = 0.0001280914293812323
$ws.Range("J5").Value2 = 0.0001280914293812324
$ws.Range("K5").Value2 = 3
$ws.Range("L5").Value2 = 1
$ws.Range("M5").Value2 = 0.5668716666666667
$ws.Range("N5").Value2 = 1.700615
$ws.Range("O5").Value2 = 0.9911343043325463
$ws.Range("P5").Value2 = 0.9911343043325463
$ws.Range("Q5").Value2 = 0.006544722348888888
$ws.Range("R5").Value2 = 0.05890250114
$ws.Range("S5").Value2 = 0.0001269558097507292
$ws.Range("T5").Value2 = 0.0001269558097507292

# Row 6
$ws.Range("A6").Value2 = "Resolving-Mac"
$ws.Range("B6").Value2 = "Cd86"
$ws.Range("C6").Value2 = "Ctla4"
$ws.Range("D6").Value2 = "MuSCs"
$ws.Range("E6").Value2 = 3
$ws.Range("F6").Value2 = 1
$ws.Range("G6").Value2 = 90.10142766666667
$ws.Range("H6").Value2 = 270.304283
$ws.Range("I6").Value2 = 0.9996437803828139
$ws.Range("J6").Value2 = 0.9996437803828139
$ws.Range("K6").Value2 = 1
$ws.Range("L6").Value2 = 0.3333333333333333
$ws.Range("M6").Value2 = 0.005070666666666667
$ws.Range("N6").Value2 = 0.015212
$ws.Range("O6").Value2 = 0.008865695667453655
$ws.Range("P6").Value2 = 0.008865695667453653
$ws.Range("Q6").Value2 = 0.4568743058884445
$ws.Range("R6").Value2 = 4.111868752996
$ws.Range("S6").Value2 = 0.008862537532736905
$ws.Range("T6").Value2 = 0.008862537532736903

# Row 7
$ws.Range("A7").Value2 = "Resolving-Mac"
$ws.Range("B7").Value2 = "Cd86"
$ws.Range("C7").Value2 = "Ctla4"
$ws.Range("D7").Value2 = "Resolving-Mac"
$ws.Range("E7").Value2 = 3
$ws.Range("F7").Value2 = 1
$ws.Range("G7").Value2 = 90.10142766666667
$ws.Range("H7").Value2 = 270.304283
$ws.Range("I7").Value2 = 0.9996437803828139
$ws.Range("J7").Value2 = 0.9996437803828139
$ws.Range("K7").Value2 = 3
$ws.Range("L7").Value2 = 1
$ws.Range("M7").Value2 = 0.5668716666666667
$ws.Range("N7").Value2 = 1.700615
$ws.Range("O7").Value2 = 0.9911343043325463
$ws.Range("P7").Value2 = 0.9911343043325463
$ws.Range("Q7").Value2 = 51.07594647044944
$ws.Range("R7").Value2 = 459.683518234045
$ws.Range("S7").Value2 = 0.990781242850077
$ws.Range("T7").Value2 = 0.990781242850077
